# "Minor update on ppt"
#  1) Bump the auto-generated "datetimeFigureOut" footer field from
#     6/15/2021 -> 6/16/2021 everywhere it appears (the slide master and
#     every slide layout each carry their own cached copy of the field).
#  2) Remove the speaker-credit subtitle ("Sarah Bauduin - ...") and the
#     accompanying headshot picture from the title slide (slide 1),
#     leaving the title and the date textbox in the corner.

$p = $ppt.ActivePresentation

$oldDate = "6/15/2021"
$newDate = "6/16/2021"

# ppPlaceholderDate
$ppPlaceholderDate = 16

function Update-DateField($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        $phType = -1
        try { $phType = $sh.PlaceholderFormat.Type } catch {}
        if ($phType -eq $ppPlaceholderDate -and $sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# Slide master footer date field.
Update-DateField $p.SlideMaster.Shapes

# Every slide layout's own footer date field.
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DateField $layouts.Item($li).Shapes
}

# Slide 1: drop the "Sarah Bauduin" subtitle and her picture.
$s1 = $p.Slides.Item(1)
for ($i = $s1.Shapes.Count; $i -ge 1; $i--) {
    $sh = $s1.Shapes.Item($i)
    if ($sh.Name -eq "Image 4") {
        $sh.Cut()
    } elseif ($sh.HasTextFrame -and $sh.TextFrame.HasText -and ($sh.TextFrame.TextRange.Text -like "Sarah*Bauduin*")) {
        $sh.Cut()
    }
}
